$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1720.481741820667
$ws.Range("C2").Value = 3119.94832158292
$ws.Range("D2").Value = 3966.582166973786
